# productPrice.xlsx update:
#  - Effect Date* (column I, row 2) moves from 2025-09-23 to 2025-09-24
#  - List / Actual List / Trade / Vat Price (columns J-M, row 2) move from 2969.0 to 7947.0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe keeps this a text value (shared string) instead of Excel
# auto-converting the date-like text into a date serial number.
$ws.Range("I2").Value = "'2025-09-24"

$ws.Range("J2").Value = 7947.0
$ws.Range("K2").Value = 7947.0
$ws.Range("L2").Value = 7947.0
$ws.Range("M2").Value = 7947.0
